# CCDB-202: The property association removed.
# Delete the ASSOCIATION/association/SLOT column (column E) from Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the whole "association" column (column E), shifting the
# following columns (old F,G -> new E,F) to the left.
$ws.Range("E1:E1048576").EntireColumn.Delete()

# Select the resulting column E, mirroring the selection left behind
# by the interactive column deletion that produced this workbook.
$ws.Range("E1:E1048576").Select()
